# Update cryptos list (price + 1h volume change) per the latest scrape.
# Source rows keep their A-column rank index; only B (coin), C (link),
# D (price) and E (volume/1h) cells are touched, matching the refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.077.96'
$ws.Range("E2").Value = '  +2.45%  '

$ws.Range("D3").Value = '3.131.23'
$ws.Range("E3").Value = '  +0.90%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").Value = '''587.49'
$ws.Range("E5").Value = '  +0.66%  '

$ws.Range("D6").Value = '''147.70'
$ws.Range("E6").Value = '  +3.33%  '

$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").Value = '3.123.82'
$ws.Range("E8").Value = '  +1.01%  '

$ws.Range("D9").Value = '''0.532'
$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("D10").Value = '''0.161'
$ws.Range("E10").Value = '  +11.94%  '

$ws.Range("D11").Value = '''5.75'
$ws.Range("E11").Value = '  +0.03%  '

$ws.Range("D12").Value = '''0.467'
$ws.Range("E12").Value = '  -0.10%  '

$ws.Range("D13").Value = '''0.0000253'
$ws.Range("E13").Value = '  +4.32%  '

$ws.Range("D14").Value = '''37.31'
$ws.Range("E14").Value = '  +5.05%  '

$ws.Range("E15").Value = '  -0.85%  '

$ws.Range("D16").Value = '3.652.90'
$ws.Range("E16").Value = '  +0.91%  '

$ws.Range("D17").Value = '63.961.86'
$ws.Range("E17").Value = '  +2.26%  '

# Rows 18/19 swapped rank order: WrappedEther now outranks Polkadot.
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.137.40'
$ws.Range("E18").Value = '  +1.04%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '''7.15'
$ws.Range("E19").Value = '  -1.70%  '

$ws.Range("D20").Value = '''466.65'
$ws.Range("E20").Value = '  +2.91%  '

$ws.Range("D21").Value = '''14.33'
$ws.Range("E21").Value = '  +1.83%  '

$ws.Range("D22").Value = '''0.731'
$ws.Range("E22").Value = '  -0.33%  '

$ws.Range("D23").Value = '''7.54'
$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").Value = '''13.22'
$ws.Range("E24").Value = '  -3.43%  '

$ws.Range("D25").Value = '''82.45'
$ws.Range("E25").Value = '  +0.49%  '

$ws.Range("D27").Value = '''8.96'
$ws.Range("E27").Value = '  +8.97%  '

$ws.Range("E28").Value = '  +0.69%  '

$ws.Range("E29").Value = '  -1.52%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").Value = '''6.85'
$ws.Range("E31").Value = '  +0.15%  '

$ws.Range("D32").Value = '''27.04'
$ws.Range("E32").Value = '  +0.01%  '

$ws.Range("E33").Value = '  -2.38%  '

$ws.Range("D34").Value = '0.0₃0887'
$ws.Range("E34").Value = '  +11.10%  '

$ws.Range("D35").Value = '''2.38'
$ws.Range("E35").Value = '  +8.31%  '

$ws.Range("E36").Value = '  +1.30%  '

$ws.Range("E37").Value = '  +13.14%  '

$ws.Range("D38").Value = '''6.08'
$ws.Range("E38").Value = '  +0.06%  '

# Rows 39/40 swapped rank order: Bittensor now outranks OKB.
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '''456.01'
$ws.Range("E39").Value = '  +7.73%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '''50.92'
$ws.Range("E40").Value = '  +0.63%  '

$ws.Range("D41").Value = '''8.70'
$ws.Range("E41").Value = '  -1.18%  '

$ws.Range("D42").Value = '''0.0373'
$ws.Range("E42").Value = '  +0.20%  '

$ws.Range("D43").Value = '2.892.75'
$ws.Range("E43").Value = '  -1.32%  '

$ws.Range("D44").Value = '''0.278'
$ws.Range("E44").Value = '  -1.46%  '

$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("D46").Value = '''2.17'
$ws.Range("E46").Value = '  +0.91%  '

$ws.Range("D47").Value = '''36.06'
$ws.Range("E47").Value = '  +3.34%  '

$ws.Range("D48").Value = '''125.41'
$ws.Range("E48").Value = '  -0.02%  '

$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("D51").Value = '''24.75'
$ws.Range("E51").Value = '  -0.10%  '
